# Update the last-checked timestamps recorded on Sheet1 (column B, rows 2-5)
# from the previous run (2023-07-18) to the new run (2023-07-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 7 -Day 23 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 5; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    if ($cell.Value -ne $null) {
        $cell.Value = $newDate
    }
}
